$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 3682.1853
$ws.Range("J17").Value = 3739.151
$ws.Range("L17").Value = 11217.453
$ws.Range("N17").Value = -11553.453
# Row 43
$ws.Range("H43").Value = 2042.8572
$ws.Range("I43").Value = 1467
$ws.Range("J43").Value = 2474.75
$ws.Range("K43").Value = 1467
$ws.Range("L43").Value = 2474.75
$ws.Range("M43").Value = -1398
$ws.Range("N43").Value = -2612.75
# Row 92
$ws.Range("H92").Value = 1305.3334
$ws.Range("I92").Value = 1093.5
$ws.Range("J92").Value = 3000
$ws.Range("K92").Value = 1093.5
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = 154.5
$ws.Range("N92").Value = -5496
# Row 107
$ws.Range("H107").Value = 1500
$ws.Range("I107").Value = 1500
$ws.Range("K107").Value = 1500
$ws.Range("M107").Value = 420
# Row 129
$ws.Range("H129").Value = 859.8570999999999
$ws.Range("J129").Value = 859.8570999999999
$ws.Range("L129").Value = 2579.5713
$ws.Range("N129").Value = -12579.5713
# Row 135
$ws.Range("H135").Value = 17863596
$ws.Range("I135").Value = 635.1739
$ws.Range("J135").Value = 100033220
$ws.Range("K135").Value = 5716.5651
$ws.Range("L135").Value = 900298980
$ws.Range("M135").Value = -3181.5651
$ws.Range("N135").Value = -900304050

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
# Row 88
$ws.Range("H88").Value = 126859.375
$ws.Range("I88").Value = 1349.75
$ws.Range("J88").Value = 252369
$ws.Range("K88").Value = 1349.75
$ws.Range("L88").Value = 252369
$ws.Range("M88").Value = -943.75
$ws.Range("N88").Value = -253181
# Row 91
$ws.Range("H91").Value = 126859.375
$ws.Range("I91").Value = 1349.75
$ws.Range("J91").Value = 252369
$ws.Range("K91").Value = 1349.75
$ws.Range("L91").Value = 252369
$ws.Range("M91").Value = 54.25
$ws.Range("N91").Value = -255177
# Row 102
$ws.Range("H102").Value = 1725.75
$ws.Range("I102").Value = 1530
$ws.Range("J102").Value = 1999.8
$ws.Range("K102").Value = 1530
$ws.Range("L102").Value = 1999.8
$ws.Range("M102").Value = 92
$ws.Range("N102").Value = -5243.8
# Row 132
$ws.Range("H132").Value = 13218.559
$ws.Range("I132").Value = 1352.4
$ws.Range("K132").Value = 4057.2
$ws.Range("M132").Value = -1527.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 5002574.5
$ws.Range("I7").Value = 5000149
$ws.Range("K7").Value = 5000149
$ws.Range("M7").Value = -5000036
# Row 20
$ws.Range("H20").Value = 3005.5293
$ws.Range("I20").Value = 4354.75
$ws.Range("J20").Value = 1806.2222
$ws.Range("K20").Value = 4354.75
$ws.Range("L20").Value = 1806.2222
$ws.Range("M20").Value = -4107.75
$ws.Range("N20").Value = -2300.2222
# Row 86
$ws.Range("H86").Value = 1470.7954
$ws.Range("I86").Value = 1367.0834
$ws.Range("J86").Value = 1937.5
$ws.Range("K86").Value = 1367.0834
$ws.Range("L86").Value = 1937.5
$ws.Range("M86").Value = -244.0834
$ws.Range("N86").Value = -4183.5
# Row 89
$ws.Range("H89").Value = 1470.7954
$ws.Range("I89").Value = 1367.0834
$ws.Range("J89").Value = 1937.5
$ws.Range("K89").Value = 6835.416999999999
$ws.Range("L89").Value = 9687.5
$ws.Range("M89").Value = -1219.416999999999
$ws.Range("N89").Value = -20919.5
# Row 94
$ws.Range("H94").Value = 849.3269
$ws.Range("I94").Value = 798.48834
$ws.Range("K94").Value = 798.48834
$ws.Range("M94").Value = -347.48834
# Row 105
$ws.Range("H105").Value = 4169250.8
$ws.Range("I105").Value = 2202
$ws.Range("J105").Value = 7145714.5
$ws.Range("K105").Value = 2202
$ws.Range("L105").Value = 7145714.5
$ws.Range("M105").Value = -455
$ws.Range("N105").Value = -7149208.5
# Row 134
$ws.Range("H134").Value = 29858.895
$ws.Range("I134").Value = 46401.668
$ws.Range("J134").Value = 1499.8572
$ws.Range("K134").Value = 139205.004
$ws.Range("L134").Value = 4499.571599999999
$ws.Range("M134").Value = -136670.004
$ws.Range("N134").Value = -9569.571599999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 5953117
$ws.Range("I105").Value = 9615882
$ws.Range("K105").Value = 9615882
$ws.Range("M105").Value = -9614135
# Row 132
$ws.Range("H132").Value = 15588.658
$ws.Range("I132").Value = 18424.768
$ws.Range("J132").Value = 4953.25
$ws.Range("K132").Value = 55274.304
$ws.Range("L132").Value = 14859.75
$ws.Range("M132").Value = -52744.304
$ws.Range("N132").Value = -19919.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 278.5
$ws.Range("I14").Value = 278.5
$ws.Range("K14").Value = 835.5
$ws.Range("M14").Value = -662.5
# Row 98
$ws.Range("H98").Value = 649.5
$ws.Range("J98").Value = 800
$ws.Range("L98").Value = 2400
$ws.Range("N98").Value = -5396
# Row 122
$ws.Range("H122").Value = 620.8889
$ws.Range("I122").Value = 273.33334
$ws.Range("K122").Value = 2460.00006
$ws.Range("M122").Value = -10.0000600000003
# Row 131
$ws.Range("H131").Value = 710.8099999999999
$ws.Range("J131").Value = 711.9293
$ws.Range("L131").Value = 2135.7879
$ws.Range("N131").Value = -12215.7879
# Row 138
$ws.Range("H138").Value = 126494.375
$ws.Range("I138").Value = 1538.1818
$ws.Range("J138").Value = 232226.53
$ws.Range("K138").Value = 4614.5454
$ws.Range("L138").Value = 696679.59
$ws.Range("M138").Value = 525.4546
$ws.Range("N138").Value = -706959.59

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1450
$ws.Range("I97").Value = 1450
$ws.Range("K97").Value = 1450
$ws.Range("M97").Value = -954
# Row 126
$ws.Range("H126").Value = 3798.658
$ws.Range("I126").Value = 2943.7083
$ws.Range("K126").Value = 8831.124899999999
$ws.Range("M126").Value = -6361.124899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3356.375
$ws.Range("I7").Value = 3325.1667
$ws.Range("K7").Value = 3325.1667
$ws.Range("M7").Value = -3213.1667
# Row 122
$ws.Range("H122").Value = 1156467.5
$ws.Range("I122").Value = 1636454
$ws.Range("J122").Value = 4499.8
$ws.Range("K122").Value = 4909362
$ws.Range("L122").Value = 13499.4
$ws.Range("M122").Value = -4906912
$ws.Range("N122").Value = -18399.4
# Row 126
$ws.Range("H126").Value = 3356.375
$ws.Range("I126").Value = 3325.1667
$ws.Range("K126").Value = 9975.500100000001
$ws.Range("M126").Value = -7505.500100000001
# Row 132
$ws.Range("H132").Value = 1723.4333
$ws.Range("I132").Value = 1105.5264
$ws.Range("J132").Value = 2790.7273
$ws.Range("K132").Value = 3316.5792
$ws.Range("L132").Value = 8372.1819
$ws.Range("M132").Value = -786.5792000000001
$ws.Range("N132").Value = -13432.1819

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 127
$ws.Range("H127").Value = 30000
$ws.Range("J127").Value = 30000
$ws.Range("L127").Value = 30000
$ws.Range("N127").Value = -39920

